# Daily attendance processing - 2025-11-15 06:53:06
#
# The "Recorded By" column (G) lists the contributors to each attendance
# record as a comma-separated string. For every row whose value is a
# two-part "A, B" combination that does NOT involve the backup account
# (backup@backdoor.com), the order of the two names is swapped (B, A).
# Rows with the backup account, or with only a single contributor, are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }

    if ($value -like "*backup@backdoor.com*") { continue }

    $parts = $value -split ", "
    if ($parts.Count -eq 2) {
        $newValue = $parts[1] + ", " + $parts[0]
        $cell.Value = $newValue
    }
}
